# Apply the latest crypto market snapshot to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string as TEXT (matches this sheet's
# existing convention of storing Price/Volume figures as inline strings)
# without leaving the cell's visual style changed.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 (Bitcoin)
$ws.Range('D2').Value = '26.666.14'
$ws.Range('E2').Value = '  +2.98%  '

# Row 3 (Ethereum)
$ws.Range('D3').Value = '1.689.57'
$ws.Range('E3').Value = '  +3.46%  '

# Row 4 (TetherUSD)
$ws.Range('E4').Value = '  -0.18%  '

# Row 5 (BNB)
Set-TextValue $ws.Range('D5') '217.85'
$ws.Range('E5').Value = '  +4.11%  '

# Row 6 (XRP)
Set-TextValue $ws.Range('D6') '0.5339'
$ws.Range('E6').Value = '  +2.58%  '

# Row 7 (USDC)
$ws.Range('E7').Value = '  -0.16%  '

# Row 8 (Cardano)
$ws.Range('E8').Value = '  +4.63%  '

# Row 9 (Dogecoin)
Set-TextValue $ws.Range('D9') '0.06434'
$ws.Range('E9').Value = '  +3.20%  '

# Row 10 (Solana)
Set-TextValue $ws.Range('D10') '21.66'
$ws.Range('E10').Value = '  +6.90%  '

# Row 11 (TRON)
Set-TextValue $ws.Range('D11') '0.07791'
$ws.Range('E11').Value = '  +3.05%  '

# Row 12 (WrappedEther)
$ws.Range('D12').Value = '1.696.26'
$ws.Range('E12').Value = '  +3.65%  '

# Row 13 (Polkadot)
Set-TextValue $ws.Range('D13') '4.505'
$ws.Range('E13').Value = '  +3.46%  '

# Row 14 (Polygon)
Set-TextValue $ws.Range('D14') '0.5629'
$ws.Range('E14').Value = '  +4.01%  '

# Row 15 (ShibaInu)
$ws.Range('D15').Value = '0.0₅8468'
$ws.Range('E15').Value = '  +6.99%  '

# Row 16 (Litecoin)
Set-TextValue $ws.Range('D16') '66.39'
$ws.Range('E16').Value = '  +2.94%  '

# Row 17 (WrappedBTC)
$ws.Range('D17').Value = '26.703.42'
$ws.Range('E17').Value = '  +3.08%  '

# Row 18 (Dai)
$ws.Range('E18').Value = '  -0.17%  '

# Row 19 (Uniswap)
Set-TextValue $ws.Range('D19') '4.812'
$ws.Range('E19').Value = '  +4.41%  '

# Row 20 (BitcoinCash)
Set-TextValue $ws.Range('D20') '195.60'
$ws.Range('E20').Value = '  +6.17%  '

# Row 21 (Avalanche)
$ws.Range('E21').Value = '  +4.20%  '

# Row 22 (Chainlink)
$ws.Range('E22').Value = '  +5.36%  '

# Row 23 (BinanceUSD)
$ws.Range('E23').Value = '  -0.27%  '

# Row 24 (Monero)
Set-TextValue $ws.Range('D24') '143.96'
$ws.Range('E24').Value = '  -1.09%  '

# Row 25 (Stellar)
Set-TextValue $ws.Range('D25') '0.1289'
$ws.Range('E25').Value = '  +7.50%  '

# Row 26 (Cosmos)
Set-TextValue $ws.Range('D26') '7.489'
$ws.Range('E26').Value = '  +2.06%  '

# Row 27 (EthereumClassic)
Set-TextValue $ws.Range('D27') '16.30'
$ws.Range('E27').Value = '  +5.29%  '

# Row 28 (Toncoin)
Set-TextValue $ws.Range('D28') '1.425'
$ws.Range('E28').Value = '  +3.75%  '

# Row 29 (Hedera)
Set-TextValue $ws.Range('D29') '0.06176'
$ws.Range('E29').Value = '  +3.87%  '

# Row 30 (PancakeSwap)
$ws.Range('E30').Value = '  +3.24%  '

# Row 31 (InternetComputer(DFINITY))
Set-TextValue $ws.Range('D31') '3.604'
$ws.Range('E31').Value = '  +7.85%  '

# Row 32 (Filecoin)
$ws.Range('E32').Value = '  +3.84%  '

# Row 34 (ARBITRUM)
Set-TextValue $ws.Range('D34') '1.015'
$ws.Range('E34').Value = '  +4.75%  '

# Row 35 (MXToken)
Set-TextValue $ws.Range('D35') '2.802'
$ws.Range('E35').Value = '  +2.41%  '

# Row 36 (HuobiToken)
$ws.Range('E36').Value = '  +1.71%  '

# Row 37 (ImmutableX)
Set-TextValue $ws.Range('D37') '0.5748'
$ws.Range('E37').Value = '  -0.93%  '

# Row 38 (VeChain)
Set-TextValue $ws.Range('D38') '0.01653'
$ws.Range('E38').Value = '  +3.95%  '

# Row 39 (FraxShare)
$ws.Range('E39').Value = '  +6.47%  '

# Row 40 (Maker)
$ws.Range('D40').Value = '1.080.64'
$ws.Range('E40').Value = '  +5.84%  '

# Row 41 (TrustWalletToken)
Set-TextValue $ws.Range('D41') '0.8674'
$ws.Range('E41').Value = '  +3.29%  '

# Row 42 (PaxDollar)
Set-TextValue $ws.Range('D42') '1.001'
$ws.Range('E42').Value = '  -0.04%  '

# Row 43 (Quant)
Set-TextValue $ws.Range('D43') '100.50'
$ws.Range('E43').Value = '  +0.97%  '

# Row 44 (RocketPoolETH)
$ws.Range('D44').Value = '1.840.45'
$ws.Range('E44').Value = '  +3.11%  '

# BabyDogeCoin was delisted: rows 45-50 each take on the coin data that
# used to sit one row below them, and the newly tracked Algorand becomes
# the last row. The row-index column A is untouched throughout.
# Row 45 (Aave)
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D45') '57.46'
$ws.Range('E45').Value = '  +5.91%  '

# Row 46 (EnergySwap)
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D46') '8.229'
$ws.Range('E46').Value = '  +3.50%  '

# Row 47 (Frax)
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range('D47') '1.003'
$ws.Range('E47').Value = '  +0.26%  '

# Row 48 (Cronos)
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D48') '0.05224'
$ws.Range('E48').Value = '  +0.98%  '

# Row 49 (Aptos)
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D49') '6.106'
$ws.Range('E49').Value = '  +5.67%  '

# Row 50 (Mantle)
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D50') '0.4240'
$ws.Range('E50').Value = '  +0.32%  '

# Row 51 (Algorand)
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D51') '0.09928'
$ws.Range('E51').Value = '  +3.82%  '

